$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = '67.549.69'
$ws.Cells.Item(2,5).Value = '  -1.67%  '

$ws.Cells.Item(3,4).Value = '3.516.56'
$ws.Cells.Item(3,5).Value = '  -4.00%  '

$ws.Cells.Item(4,5).Value = '  +0.18%  '

$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '198.85'
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = '  -0.51%  '

$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '552.51'
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = '  -4.73%  '

$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = '0.633'
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).Value = '  +1.80%  '

$ws.Cells.Item(8,4).Value = '3.506.43'
$ws.Cells.Item(8,5).Value = '  -4.13%  '

$ws.Cells.Item(9,5).Value = '  +0.01%  '

$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '0.654'
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value = '  -4.14%  '

$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '61.85'
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = '  +8.31%  '

$ws.Cells.Item(12,5).Value = '  -8.05%  '

$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = '0.0000269'
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).Value = '  -9.53%  '

$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = '9.83'
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value = '  -3.36%  '

$ws.Cells.Item(15,4).Value = '4.070.19'
$ws.Cells.Item(15,5).Value = '  -4.15%  '

$ws.Cells.Item(16,4).Value = '3.513.10'
$ws.Cells.Item(16,5).Value = '  -3.82%  '

$ws.Cells.Item(17,5).Value = '  -1.96%  '

$ws.Cells.Item(18,4).Value = '67.275.20'
$ws.Cells.Item(18,5).Value = '  -1.89%  '

$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = '18.36'
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Value = '  -1.95%  '

$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = '11.82'
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value = '  -6.28%  '

$ws.Cells.Item(21,5).Value = '  -6.38%  '

$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = '393.81'
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value = '  -3.11%  '

$ws.Cells.Item(23,5).Value = '  -7.07%  '

$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = '11.87'
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).Value = '  -9.92%  '

$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = '85.00'
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value = '  -1.50%  '

$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '3.87'
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).Value = '  -0.58%  '

$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = '12.27'
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value = '  -3.42%  '

$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = '2.82'
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).Value = '  -5.36%  '

$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = '8.85'
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).Value = '  -4.71%  '

$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = '719.12'
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).Value = '  +3.26%  '

$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = '31.10'
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).Value = '  -3.36%  '

$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = '7.02'
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,5).Value = '  -15.64%  '

$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = '11.71'
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).Value = '  -5.15%  '

$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = '63.96'
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).Value = '  -1.39%  '

$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = '0.111'
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).Value = '  -5.82%  '

$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = '38.50'
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).Value = '  -10.69%  '

$ws.Cells.Item(37,5).Value = '  -0.20%  '

$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = '0.394'
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).Value = '  -8.25%  '

$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = '3.02'
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).Value = '  -4.55%  '

$ws.Cells.Item(40,2).Value = 'Kaspa'
$ws.Cells.Item(40,3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = '0.130'
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value = '  -8.68%  '

$ws.Cells.Item(41,2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(41,3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = '0.998'
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value = '  +0.05%  '

$ws.Cells.Item(42,4).Value = '3.066.53'
$ws.Cells.Item(42,5).Value = '  -4.90%  '

$ws.Cells.Item(43,5).Value = '  -15.45%  '

$ws.Cells.Item(44,5).Value = '  +5.57%  '

$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = '2.51'
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value = '  -13.66%  '

$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = '0.0407'
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Value = '  -4.11%  '

$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = '0.130'
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).Value = '  -2.02%  '

$ws.Cells.Item(48,5).Value = '  -14.10%  '

$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = '138.56'
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Value = '  -3.06%  '

$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = '8.24'
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value = '  -8.38%  '

$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = '2.90'
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Value = '  -6.60%  '
